$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the bet date (A1): 45755 -> 45756
$ws.Range("A1").Value = 45756

# --- First parlay leg group (rows 2-5) ---
# A2:A4 share the same string ("ORL vs ATL" -> "ORL vs BOS")
$ws.Range("A2").Value = "ORL vs BOS"
$ws.Range("A3").Value = "ORL vs BOS"
$ws.Range("A4").Value = "ORL vs BOS"

$ws.Range("B2").Value = "Kentavious Caldwell-Pope - Sobre 2.5 Goles de campo realizados"
$ws.Range("C2").Value = "1,54x"
$ws.Range("D2").Value = "Cumplio en 8/10 https://prnt.sc/EPEzM4XtY5Io"

$ws.Range("B3").Value = "Neemias Queta - Sobre 1.5 Faltas personales"
$ws.Range("C3").Value = "1,47x"
$ws.Range("D3").Value = "Cumplio en 10/10 sin Horford ni Porzingis (pivot suplente de kornet hoy) https://prnt.sc/1r1GPmXWdrfG"

$ws.Range("B4").Value = "Sam Hauser - Sobre 4.5 Goles de campo realizados"
$ws.Range("C4").Value = "1,53x"
$ws.Range("D4").Value = "Cumplio en 2/2 (contando los out que hay Tatum, porzingis white, holiday, horford)https://prnt.sc/Cy_6ypnFUpPy"

$ws.Range("C5").Value = "4,00x"

# --- Second parlay leg group (rows 6-9) ---
$ws.Range("A6").Value = "PHI vs WAS"
$ws.Range("B6").Value = "Jared Butler - Sobre 1.5 Faltas personales"
$ws.Range("C6").Value = "1,62x"
$ws.Range("D6").Value = "Cumplio en 8/10 https://prnt.sc/XTQsqZhitf3W"

$ws.Range("A7").Value = "PHI vs WAS"
$ws.Range("B7").Value = "Quentin Grimes - Sobre 22.5 Puntos"
$ws.Range("C7").Value = "1,47x"
$ws.Range("D7").Value = "Cumplio en 8/10 https://prnt.sc/gV7LSMhZRTKG"

$ws.Range("A8").Value = "PHI vs WAS"
$ws.Range("B8").Value = "Kyshawn George - Sobre 3.5 Rebotes"
$ws.Range("C8").Value = "1,53x"
$ws.Range("D8").Value = "Cumplio en 7/10 (hizo 3) https://prnt.sc/lUx3GAQ6NLxo"

$ws.Range("C9").Value = "4,20x"

# --- Third parlay leg group (rows 10-13) ---
$ws.Range("A10").Value = "CHI vs MIA"
$ws.Range("B10").Value = "Nikola Vucevic - Sobre 7.5 Rebotes"
$ws.Range("C10").Value = "1,35x"
$ws.Range("D10").Value = "Cumplio en 9/10 https://prnt.sc/mD4dP-gNddaU"

$ws.Range("A11").Value = "CHI vs MIA"
$ws.Range("B11").Value = "Davion Mitchell - Sobre 7.5 Puntos"
$ws.Range("C11").Value = "1,46x"
$ws.Range("D11").Value = "Cumplio en 9/10 (hizo 6) https://prnt.sc/Dn1Y-sz7MTHt"

$ws.Range("A12").Value = "CHI vs MIA"
$ws.Range("B12").Value = "Coby White - Sobre 28.5 Puntos + Asistencias + Rebotes"
$ws.Range("C12").Value = "1,45x"
$ws.Range("D12").Value = "Cumplio en 9/10 https://prnt.sc/cCpSsKkGug-V"

$ws.Range("C13").Value = "3,15x"

# --- Totals row ---
$ws.Range("C14").Value = "52,26x"
$ws.Range("D14").Value = "https://stake.com/sports/home?betId=cad1a1e1-360a-4cb2-bae4-f3b02ec8f0db&modal=bet"

# Row heights for the taller info rows (3 and 4)
$ws.Rows.Item(3).RowHeight = 35.25
$ws.Rows.Item(4).RowHeight = 36

# Selection moves to B4
$ws.Range("B4").Select()
